$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 2.18
$ws.Range("H7").Value = 2.87
$ws.Range("I7").Value = 3.55
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.45
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 6.3
$ws.Range("Z7").Value = 22
$ws.Range("AG7").Value = 9
$ws.Range("AH7").Value = 18.5
$ws.Range("AM7").Value = 700
$ws.Range("AU7").Value = 6.7
$ws.Range("AV7").Value = 60
$ws.Range("AY7").Value = 25
$ws.Range("AZ7").Value = 100
$ws.Range("BB7").Value = 300
$ws.Range("G8").Value = 2.5
$ws.Range("I8").Value = 2.75
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 3.5
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65
$ws.Range("X8").Value = 12
$ws.Range("Z8").Value = 26
$ws.Range("AK8").Value = 23
$ws.Range("AR8").Value = 81
$ws.Range("AY8").Value = 26
$ws.Range("G10").Value = 2.35
$ws.Range("H10").Value = 2.65
$ws.Range("K10").Value = 1.87
$ws.Range("Q11").Value = 1.89
$ws.Range("R11").Value = 1.84
$ws.Range("G12").Value = 1.22
$ws.Range("H12").Value = 5.25
$ws.Range("I12").Value = 13
$ws.Range("J12").Value = 1.69
$ws.Range("L12").Value = 11
$ws.Range("Q12").Value = 1.87
$ws.Range("R12").Value = 1.87
$ws.Range("X12").Value = 5
$ws.Range("AD12").Value = 11
$ws.Range("AE12").Value = 34
$ws.Range("AF12").Value = 126
$ws.Range("AP12").Value = 23
$ws.Range("AQ12").Value = 15
$ws.Range("AV12").Value = 101
$ws.Range("I14").Value = 2.47
$ws.Range("P14").Value = 2.65
$ws.Range("V14").Value = 1.83
$ws.Range("AF14").Value = 70
$ws.Range("AG14").Value = 7.2
$ws.Range("AL14").Value = 32
